# Wishlist workbook update: add a "Reserved / bought" tracking column,
# mark the items that have been reserved/bought, turn the table into an
# AutoFilter range, and leave the cursor on the newly added column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the old "Reserved" header (column E) to "Reserved / bought"
$ws.Range("E1").Value = "Reserved / bought"

# Mark the rows that have been reserved / bought with a "Y"
$reservedRows = @(5, 7, 8, 9, 10, 15, 16)
foreach ($r in $reservedRows) {
    $ws.Cells.Item($r, 5).Value = "Y"
}

# Size the new column to fit its contents
$ws.Columns("E").ColumnWidth = 17.5703125

# Turn the table (A1:E16) into a filterable range
$ws.Range("A1:E16").AutoFilter()

# Register the (hidden) sheet-scoped _FilterDatabase name that Excel
# creates behind the scenes whenever AutoFilter is turned on
$filterDatabase = $ws.Names.Add("_xlnm._FilterDatabase", "=Raluca!`$A`$1:`$E`$16")
$filterDatabase.Visible = $false

# Leave the selection where the user last clicked
$ws.Range("E13").Select()
